$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new coefficient values for row 4 (id = 3)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "0.5"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "0.5"
$ws.Range("I4").Value = 0

# Move the active selection to B5 as left by the editor
$ws.Range("B5").Select()
